$d = $word.ActiveDocument

# Merge the split-run Title paragraph into a single run.
$d.Content.Find.Execute(
    "Questions: Introduction to complex numbers",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Questions: Introduction to complex numbers", 2) | Out-Null

# Merge the split-run Author paragraph into a single run.
$d.Content.Find.Execute(
    "Tom Coleman",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Tom Coleman", 2) | Out-Null

# Merge the split-run Abstract paragraph into a single run.
$d.Content.Find.Execute(
    "A selection of questions for the study guide on introduction to complex numbers.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "A selection of questions for the study guide on introduction to complex numbers.", 2) | Out-Null
